$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the text in F21 (row 21 = displayed "22"), which holds the long text
# about MasterPage limitation.
$ws.Range("F21").Value = "Không làm được với MasterPage và web user control"

# Widen column F to fit the new text (maps to stored OOXML width of 41).
$ws.Columns("F").ColumnWidth = 40.140625

# Move the active selection/cursor to A20 (matches the author's recorded cursor position).
$ws.Range("A20").Select()
